# Updates the cryptos price/volume table (columns D and E) to reflect
# the latest GitHub Actions refresh of the data, per the commit:
# "Updated cryptos list on Sun Apr 16 03:27:06 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.563.82"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.109.04"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D5").Value = "336.09"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("D7").Value = "0.5247"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("D8").Value = "0.4541"
$ws.Range("E8").Value = "  +4.31%  "
$ws.Range("D9").Value = "55.80"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("D10").Value = "0.09041"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").Value = "1.168"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "24.56"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "2.119.13"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "6.846"
$ws.Range("E14").Value = "  +2.53%  "
$ws.Range("D15").Value = "8.107"
$ws.Range("E15").Value = "  +5.56%  "
$ws.Range("D16").Value = "0.00001180"
$ws.Range("E16").Value = "  +5.36%  "
$ws.Range("D17").Value = "97.20"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "1.012"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "0.06690"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "1.010"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "6.259"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "30.625.56"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").Value = "12.77"
$ws.Range("E24").Value = "  +4.46%  "
$ws.Range("D25").Value = "2.362"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").Value = "2.364.66"
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "163.36"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "2.515"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "133.53"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "1.217"
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("D32").Value = "0.1068"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "6.339"
$ws.Range("E33").Value = "  +3.48%  "
$ws.Range("D34").Value = "1.619"
$ws.Range("E34").Value = "  -2.73%  "
$ws.Range("D35").Value = "3.969"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").Value = "10.45"
$ws.Range("E36").Value = "  +3.76%  "
$ws.Range("D37").Value = "5.862"
$ws.Range("E37").Value = "  +7.60%  "
$ws.Range("D38").Value = "0.02611"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").Value = "0.06817"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "0.2310"
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("D41").Value = "12.58"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("D42").Value = "0.6835"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").Value = "1.257"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "0.6437"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D45").Value = "14.02"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").Value = "2.303"
$ws.Range("E46").Value = "  +4.75%  "
$ws.Range("D47").Value = "3.673"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").Value = "1.250"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "0.00000000348"
$ws.Range("E49").Value = "  +17.91%  "
$ws.Range("D50").Value = "1.211"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").Value = "83.13"
$ws.Range("E51").Value = "  +1.55%  "
